# Auto-generated edit script: update crypto price/volume table
# to reflect the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.324.24"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.932.37"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'325.34"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4621"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'0.3871"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "'45.91"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'0.07827"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").Value = "'0.9737"
$ws.Range("E11").Value = "  -1.79%  "
$ws.Range("D12").Value = "'22.59"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").Value = "1.922.43"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "'7.076"
$ws.Range("D15").Value = "'5.764"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "'0.07051"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'86.66"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "'0.000009794"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("D20").Value = "'17.05"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "29.396.44"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").Value = "'0.4969"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "'5.482"
$ws.Range("E24").Value = "  +3.12%  "
$ws.Range("D26").Value = "2.168.27"
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("D28").Value = "'157.70"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").Value = "'19.40"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'5.744"
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("D31").Value = "'118.39"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'1.856"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").Value = "'0.09352"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "'0.8588"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "'5.171"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").Value = "'1.305"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "'3.087"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("D38").Value = "'0.05769"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").Value = "'1.154"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "'0.02075"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Value = "'7.682"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "'0.5660"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").Value = "'0.1777"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "'0.000002997"
$ws.Range("E44").Value = "  +51.03%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'9.402"
$ws.Range("E45").Value = "  -3.25%  "
$ws.Range("D46").Value = "'2.692"
$ws.Range("E46").Value = "  +5.47%  "
$ws.Range("D47").Value = "'0.5283"
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("D48").Value = "'11.49"
$ws.Range("E48").Value = "  -3.62%  "
$ws.Range("D49").Value = "'0.06871"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").Value = "'2.074"
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("D51").Value = "'1.816"
$ws.Range("E51").Value = "  -1.83%  "
